$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the affected cells so values round-trip as literal strings
# (matching the inlineStr / text representation used in the source data), then set values.

$targets = @("D2","E2","D3","E3","D4","E4","D5","E5","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","D17","E17","E18","D19","E19","E21","D22","E22","E23","D24","E24","D25","E25","E26","D27","E27","D39","E39","D40","E40","D41","E41","D42","E42","E43","D44","E44","D45","E45","D46","E46","E47")
foreach ($addr in $targets) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "298.63"
$ws.Range("E2").Value = "-2.33%"
$ws.Range("D3").Value = "31.76"
$ws.Range("E3").Value = "-1.46%"
$ws.Range("D4").Value = "5.107"
$ws.Range("E4").Value = "-4.17%"
$ws.Range("D5").Value = "0.07539"
$ws.Range("E5").Value = "1.78%"
$ws.Range("E6").Value = "0.05%"
$ws.Range("D7").Value = "1.752"
$ws.Range("E7").Value = "13.86%"
$ws.Range("D8").Value = "3.794"
$ws.Range("E8").Value = "2.64%"
$ws.Range("D9").Value = "0.9279"
$ws.Range("E9").Value = "1.86%"
$ws.Range("D10").Value = "0.1711"
$ws.Range("E10").Value = "2.75%"
$ws.Range("D11").Value = "0.07303"
$ws.Range("E11").Value = "-4.39%"
$ws.Range("D12").Value = "0.07934"
$ws.Range("E12").Value = "-0.17%"
$ws.Range("D13").Value = "0.03055"
$ws.Range("E13").Value = "1.02%"
$ws.Range("D14").Value = "0.09887"
$ws.Range("E14").Value = "0.32%"
$ws.Range("D15").Value = "0.001501"
$ws.Range("E15").Value = "-1.59%"
$ws.Range("D16").Value = "0.04647"
$ws.Range("D17").Value = "0.006414"
$ws.Range("E17").Value = "0.80%"
$ws.Range("E18").Value = "-1.12%"
$ws.Range("D19").Value = "2.221"
$ws.Range("E19").Value = "-0.81%"
$ws.Range("E21").Value = "-1.50%"
$ws.Range("D22").Value = "4.569"
$ws.Range("E22").Value = "8.22%"
$ws.Range("E23").Value = "-4.95%"
$ws.Range("D24").Value = "0.001215"
$ws.Range("E24").Value = "-0.13%"
$ws.Range("D25").Value = "0.004424"
$ws.Range("E25").Value = "-1.77%"
$ws.Range("E26").Value = "19.53%"
$ws.Range("D27").Value = "0.0001858"
$ws.Range("E27").Value = "6.69%"
$ws.Range("D39").Value = "0.01680"
$ws.Range("E39").Value = "1.87%"
$ws.Range("D40").Value = "0.04543"
$ws.Range("E40").Value = "1.17%"
$ws.Range("D41").Value = "0.007060"
$ws.Range("E41").Value = "-4.30%"
$ws.Range("D42").Value = "0.1328"
$ws.Range("E42").Value = "-2.31%"
$ws.Range("E43").Value = "-8.95%"
$ws.Range("D44").Value = "0.01279"
$ws.Range("E44").Value = "-6.17%"
$ws.Range("D45").Value = "0.00006003"
$ws.Range("E45").Value = "-2.25%"
$ws.Range("D46").Value = "1.930"
$ws.Range("E46").Value = "1.97%"
$ws.Range("E47").Value = "-0.26%"
